$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Strings" class (row 25) pointed its "Python Strings" pre-class
# reading at the docs.python.org tutorial; swap it for the YouTube video
# (second bullet, about Unicode, is unchanged).
$ws.Range("D25").Value = "'- ``Python Strings <https://www.youtube.com/watch?v=pUbfDilfutE>```_`n- ``Computerphile Unicode <https://www.youtube.com/watch?v=MijmeoH9LT4>```_"

# "Pandas: Indices & Missing" class (row 11) was missing its in-class
# exercise link -- add it.
$ws.Range("E11").Value = "'``Link <exercises/Exercise_indices_missing.ipynb>```_"

# Reflect the resulting selection (the cell that was just filled in).
$ws.Range("E11").Select()
